$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 38

$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/01/2026"
$ws.Cells.Item($row, 2).Value = 11770.23
$ws.Cells.Item($row, 3).Value = 0.2189261117235364
$ws.Cells.Item($row, 4).Value = 0.7810738882764636
$ws.Cells.Item($row, 5).Value = -155.88
$ws.Cells.Item($row, 6).Value = -29.14
$ws.Cells.Item($row, 7).Value = -21398.34
$ws.Cells.Item($row, 8).Value = -69.95
$ws.Cells.Item($row, 9).Value = -476.04
$ws.Cells.Item($row, 10).Value = -15.59
